$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "0.6"
$ws.Range("G3").Value = "0.54"
$ws.Range("H3").Value = "0.49"
$ws.Range("I3").Value = "0.57"
$ws.Range("H4").Value = "0.31"
$ws.Range("F5").Value = "0.12"
$ws.Range("G5").Value = "0.13"
$ws.Range("H5").Value = "0.16"
$ws.Range("I5").Value = "0.13"
$ws.Range("F6").Value = "0.03"
$ws.Range("H6").Value = "0.05"
$ws.Range("I6").Value = "0.04"
$ws.Range("F8").Value = "0.26"
$ws.Range("G8").Value = "0.09"
$ws.Range("H8").Value = "0.1"
$ws.Range("B9").Value = "0.31"
$ws.Range("C9").Value = "0.24"
$ws.Range("D9").Value = "0.27"
$ws.Range("E9").Value = "0.19"
$ws.Range("F9").Value = "0.22"
$ws.Range("G9").Value = "0.11"
$ws.Range("H9").Value = "0.1"
$ws.Range("I9").Value = "0.16"
$ws.Range("J9").Value = "0.14"
$ws.Range("K9").Value = "0.06"
$ws.Range("L9").Value = "0.09"
$ws.Range("B10").Value = "0.26"
$ws.Range("C10").Value = "0.45"
$ws.Range("D10").Value = "0.34"
$ws.Range("E10").Value = "0.48"
$ws.Range("F10").Value = "0.4"
$ws.Range("G10").Value = "0.55"
$ws.Range("H10").Value = "0.48"
$ws.Range("I10").Value = "0.57"
$ws.Range("J10").Value = "0.61"
$ws.Range("B11").Value = "0.05"
$ws.Range("C11").Value = "0.13"
$ws.Range("D11").Value = "0.09"
$ws.Range("E11").Value = "0.2"
$ws.Range("F11").Value = "0.09"
$ws.Range("G11").Value = "0.21"
$ws.Range("H11").Value = "0.3"
$ws.Range("I11").Value = "0.14"
$ws.Range("J11").Value = "0.19"
$ws.Range("K11").Value = "0.28"
$ws.Range("L11").Value = "0.27"
$ws.Range("F13").Value = "3.1"
$ws.Range("G13").Value = "2.68"
$ws.Range("H13").Value = "2.58"
$ws.Range("I13").Value = "2.36"
$ws.Range("F14").Value = "0.19"
$ws.Range("G14").Value = "0.24"
$ws.Range("H14").Value = "0.25"
$ws.Range("I14").Value = "0.3"
$ws.Range("F15").Value = "0.49"
$ws.Range("G15").Value = "0.41"
$ws.Range("I15").Value = "0.3"
$ws.Range("B16").Value = "0.5"
$ws.Range("C16").Value = "0.34"
$ws.Range("D16").Value = "0.35"
$ws.Range("E16").Value = "0.36"
$ws.Range("F16").Value = "0.41"
$ws.Range("G16").Value = "0.44"
$ws.Range("H16").Value = "0.52"
$ws.Range("I16").Value = "0.44"
$ws.Range("J16").Value = "0.29"
$ws.Range("K16").Value = "0.52"
$ws.Range("L16").Value = "0.43"
$ws.Range("G18").Value = "0.06"
$ws.Range("H18").Value = "0.12"
$ws.Range("A19").Value = "15 - 24"
$ws.Range("G19").Value = "0.19"
$ws.Range("H19").Value = "0.21"
$ws.Range("A20").Value = "25 - 49"
$ws.Range("F20").Value = "0.54"
$ws.Range("G20").Value = "0.5"
$ws.Range("H20").Value = "0.46"
$ws.Range("F21").Value = "0.17"
$ws.Range("G21").Value = "0.24"
$ws.Range("H21").Value = "0.21"
$ws.Range("G24").Value = "0.33"
$ws.Range("F25").Value = "0.44"
$ws.Range("H25").Value = "0.34"
$ws.Range("F26").Value = "0.15"
$ws.Range("G26").Value = "0.24"
$ws.Range("F28").Value = "0.68"
$ws.Range("G28").Value = "0.72"
$ws.Range("F29").Value = "0.79"
$ws.Range("G29").Value = "0.85"
$ws.Range("H29").Value = "0.91"
$ws.Range("F30").Value = "12092"
$ws.Range("G30").Value = "14473"
$ws.Range("H30").Value = "197822"
$ws.Range("I30").Value = "66362"
